$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stock List")

# Shift rows 2-75 down to rows 3-76 (columns B, C, D, E, H) so a new
# row can be inserted at row 2 without expanding the used range.
for ($r = 76; $r -ge 3; $r--) {
    $src = $r - 1
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($src, 2).Value()
    $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($src, 3).Value()
    $ws.Cells.Item($r, 4).Value = $ws.Cells.Item($src, 4).Value()
    $ws.Cells.Item($r, 5).Value = $ws.Cells.Item($src, 5).Value()
    $ws.Cells.Item($r, 8).Value = $ws.Cells.Item($src, 8).Value()
}

# Write the new first data row.
$ws.Cells.Item(2, 2).Value = "CAPTRU-RE1"
$ws.Cells.Item(2, 3).Value = "CAPTRU-RE1"
$ws.Cells.Item(2, 4).Value = 5.67
$ws.Cells.Item(2, 5).Value = -11.9565
